$d = $word.ActiveDocument

# 1) Remove the "Meta description: ..." paragraph that currently sits
#    right under the H1 title.
$metaRng = $d.Content
$metaRng.Find.Execute("Meta description", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$metaRng.Expand(4)   # wdParagraph
$metaRng.Delete()

# 2) Insert a new bold paragraph "Play Diamond Queen Free Slot by IGT"
#    right before the "Prompt for DALLE..." paragraph near the end.
$dalleRng = $d.Content
$dalleRng.Find.Execute("Prompt for DALLE", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$dalleRng.Expand(4)  # wdParagraph

$precedingRange = $d.Range(0, $dalleRng.Start)
$prevPara = $precedingRange.Paragraphs.Last
$insertAt = $prevPara.Range.End
$newTitle = "Play Diamond Queen Free Slot by IGT"

$insertPoint = $d.Range($insertAt, $insertAt)
$insertPoint.InsertAfter("$newTitle`r")

$titleRange = $d.Range($insertAt, $insertAt + $newTitle.Length)
$titleRange.Font.Bold = $true

# 3) Swap the (still italic) "Prompt for DALLE..." copy for the new
#    review blurb text.
$oldBlurb = 'Prompt for DALLE: Create a feature image for the online slot game "Diamond Queen" featuring a happy Maya warrior with glasses in a cartoon style. The image should be vibrant and eye-catching, with the Maya warrior holding a large diamond scepter, surrounded by precious jewels and enchanted forest elements. The image should convey a sense of magic and fantasy, while also highlighting the diamond theme of the game. The Maya warrior should be depicted with a joyful expression and be wearing glasses, emphasizing the technological aspect of the game. Overall, the image should appeal to players who are looking for an exciting and magical gaming experience combined with cutting-edge technology.'
$newBlurb = 'Read our review of Diamond Queen, a magical and elegant online slot by IGT. Play for free and trigger the Mystical Diamond Bonus for extra Wilds and free spins.'
$d.Content.Find.Execute($oldBlurb, $true, $false, $false, $false, $false, $true, 1, $false, $newBlurb, 2)
